$d = $word.ActiveDocument

$replacements = @(
    @{old="743×5=3715"; new="275×6=1650"},
    @{old="541×7=3787"; new="543×3=1629"},
    @{old="944×7=6608"; new="724×9=6516"},
    @{old="892×8=7136"; new="847×2=1694"},
    @{old="190×2=380";  new="439×4=1756"},
    @{old="721×9=6489"; new="173×9=1557"},
    @{old="843×9=7587"; new="417×5=2085"},
    @{old="331×7=2317"; new="318×6=1908"},
    @{old="365×2=730";  new="806×4=3224"},
    @{old="127×2=254";  new="984×6=5904"},
    @{old="982×2=1964"; new="320×6=1920"},
    @{old="535×3=1605"; new="886×6=5316"},
    @{old="616×4=2464"; new="388×4=1552"},
    @{old="612×2=1224"; new="378×4=1512"},
    @{old="635×3=1905"; new="144×2=288"},
    @{old="847×4=3388"; new="595×8=4760"},
    @{old="388×6=2328"; new="889×6=5334"},
    @{old="848×5=4240"; new="845×6=5070"},
    @{old="251×5=1255"; new="551×5=2755"},
    @{old="930×9=8370"; new="494×7=3458"},
    @{old="701×6=4206"; new="110×6=660"},
    @{old="497×7=3479"; new="437×5=2185"},
    @{old="262×7=1834"; new="598×3=1794"},
    @{old="500×7=3500"; new="198×9=1782"},
    @{old="834×3=2502"; new="144×3=432"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
